$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 59.895130354105632
$ws.Range("C2").Value = 48.282181340544462
$ws.Range("D2").Value = 62.786754089003175
$ws.Range("E2").Value = 51.396987585438829

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 61.950761800214487
$ws.Range("C3").Value = 44.044568128126528
$ws.Range("D3").Value = 71.204113636332352
$ws.Range("E3").Value = 47.21481631296556

# Update the selection to match B1:E3
$ws.Range("B1:E3").Select()
